$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("D2").Value = "43.000.02"
$ws.Range("E2").Value = "  +0.54%  "
$ws.Range("D3").Value = "2.299.70"
$ws.Range("E3").Value = "  +2.00%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").Value = "'252.81"
$ws.Range("E5").Value = "  -0.41%  "
$ws.Range("D6").Value = "'0.643"
$ws.Range("E6").Value = "  +2.69%  "
$ws.Range("D7").Value = "'74.35"
$ws.Range("E7").Value = "  +3.55%  "
$ws.Range("E8").Value = "  +0.07%  "
$ws.Range("E9").Value = "  -1.26%  "
$ws.Range("D10").Value = "'39.56"
$ws.Range("E10").Value = "  -3.31%  "
$ws.Range("D11").Value = "'0.0985"
$ws.Range("E11").Value = "  +1.52%  "
$ws.Range("D12").Value = "'7.51"
$ws.Range("E12").Value = "  +1.15%  "
$ws.Range("E13").Value = "  +2.33%  "
$ws.Range("D14").Value = "2.645.42"
$ws.Range("E14").Value = "  +1.98%  "
$ws.Range("D15").Value = "'15.39"
$ws.Range("E15").Value = "  +3.97%  "
$ws.Range("D16").Value = "'0.876"
$ws.Range("E16").Value = "  -0.95%  "
$ws.Range("D17").Value = "2.303.12"
$ws.Range("E17").Value = "  +2.25%  "
$ws.Range("D18").Value = "42.932.17"
$ws.Range("E18").Value = "  +0.48%  "
$ws.Range("E19").Value = "  +3.71%  "
$ws.Range("D20").Value = "'6.31"
$ws.Range("E20").Value = "  +0.98%  "
$ws.Range("E21").Value = "  -0.81%  "
$ws.Range("D22").Value = "'238.45"
$ws.Range("E22").Value = "  +0.82%  "
$ws.Range("E23").Value = "  +8.12%  "
$ws.Range("D24").Value = "'3.90"
$ws.Range("E24").Value = "  -1.19%  "
$ws.Range("D25").Value = "'11.65"
$ws.Range("E25").Value = "  -0.51%  "
$ws.Range("E26").Value = "  -0.02%  "
$ws.Range("D27").Value = "'2.43"
$ws.Range("E27").Value = "  -1.05%  "
$ws.Range("E28").Value = "  -1.28%  "
$ws.Range("E29").Value = "  -0.52%  "
$ws.Range("D30").Value = "'167.74"
$ws.Range("E30").Value = "  -0.05%  "
$ws.Range("E31").Value = "  +0.33%  "
$ws.Range("D32").Value = "'6.31"
$ws.Range("E32").Value = "  +2.38%  "
$ws.Range("D33").Value = "'0.0841"
$ws.Range("E33").Value = "  +6.81%  "
$ws.Range("E34").Value = "  -2.25%  "
$ws.Range("D35").Value = "'30.97"
$ws.Range("E35").Value = "  +7.23%  "
$ws.Range("E36").Value = "  +1.75%  "
$ws.Range("D37").Value = "'4.61"
$ws.Range("E37").Value = "  +11.16%  "
$ws.Range("E38").Value = "  +2.13%  "
$ws.Range("D39").Value = "'0.0310"
$ws.Range("E39").Value = "  -2.91%  "
$ws.Range("D40").Value = "'13.80"
$ws.Range("E40").Value = "  +9.72%  "
$ws.Range("D41").Value = "'2.37"
$ws.Range("E41").Value = "  +3.37%  "
$ws.Range("E42").Value = "  +0.45%  "
$ws.Range("D43").Value = "'0.219"
$ws.Range("E43").Value = "  +8.43%  "
$ws.Range("E44").Value = "  +2.23%  "
$ws.Range("D45").Value = "'62.28"
$ws.Range("E45").Value = "  -2.62%  "
$ws.Range("E46").Value = "  -2.96%  "
$ws.Range("E47").Value = "  +1.52%  "
$ws.Range("D48").Value = "'104.05"
$ws.Range("E48").Value = "  +9.86%  "
$ws.Range("E49").Value = "  -0.03%  "
$ws.Range("E50").Value = "  +0.19%  "
$ws.Range("E51").Value = "  -0.30%  "
